# Insert a new data row before the current row 283 (shifts existing rows
# 283..292 down to 284..293) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(283).Insert()

$ws.Cells.Item(283, 1).Value = 4
$ws.Cells.Item(283, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(283, 3).Value = "Los Lagos"
$ws.Cells.Item(283, 4).Value = 44747
$ws.Cells.Item(283, 5).Value = 10
$ws.Cells.Item(283, 6).Value = 100112037
$ws.Cells.Item(283, 7).Value = "Cebollín"
$ws.Cells.Item(283, 8).Value = "Sin especificar"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 160
$ws.Cells.Item(283, 11).Value = 10000
$ws.Cells.Item(283, 12).Value = 10000
$ws.Cells.Item(283, 13).Value = 10000
$ws.Cells.Item(283, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(283, 15).Value = "Región Metropolitana"
$ws.Cells.Item(283, 16).Value = 278
$ws.Cells.Item(283, 17).Value = 36
$ws.Cells.Item(283, 18).Value = "Hortaliza"
